# BOM diff tool changes:
# - The "BPP SKU" column (C) on sheet "file_1" used to just say "X" for every
#   row; it is now filled in with the same value as the "MFR PART #" column (D),
#   matching the formatting (style) already used by column D.
# - A couple of placeholder / garbage values are replaced with real-looking
#   sample text ("lorem", "eps", "Y").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("file_1")

# --- Row 4: special case, C4 becomes "lorem", D4 becomes "eps" ------------
$ws.Range("C4").Value = "lorem"
$ws.Range("D4").Value = "eps"

# --- Row 6: C6 used to hold the long garbage string; now mirrors D6 ("C") -
$ws.Range("C6").Value = "C"

# --- Row 12: special case, both C12 and D12 become "Y" --------------------
$ws.Range("C12").Value = "Y"
$ws.Range("D12").Value = "Y"

# --- Remaining rows: column C simply mirrors column D ----------------------
$ws.Range("C5").Value  = "B"
$ws.Range("C9").Value  = "A"
$ws.Range("C10").Value = "B"
$ws.Range("C11").Value = "C"
$ws.Range("C13").Value = "E"
$ws.Range("C14").Value = "F"
$ws.Range("C17").Value = "A"
$ws.Range("C18").Value = "D"
$ws.Range("C19").Value = "E"
$ws.Range("C20").Value = "F"
$ws.Range("C21").Value = "G"
$ws.Range("C22").Value = "H"
$ws.Range("C25").Value = "C"
$ws.Range("C26").Value = "F"
$ws.Range("C27").Value = "G"

# --- Match column C's formatting to column D's formatting on every data row
$dataRows = @(4,5,6,9,10,11,12,13,14,17,18,19,20,21,22,25,26,27)
foreach ($r in $dataRows) {
    $ws.Range("D$r").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- Rows 4 and 6 no longer need an explicit (taller) row height now that
# the long placeholder text is gone -----------------------------------------
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(6).EntireRow.AutoFit()

# --- Update selection to match the author's final cursor position ----------
$ws.Range("D12").Select()
